$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = New-Object 'object[,]' 24,13
$values[0,0] = 10.36439874374152
$values[0,1] = 9.08488887294193
$values[0,2] = 13.67815393050197
$values[0,3] = 33.12943604819073
$values[0,4] = 34.8586548190487
$values[0,5] = 15.91232430985558
$values[0,6] = 25.91559207776746
$values[0,7] = 10.2031243511138
$values[0,8] = 13.75449887387198
$values[0,9] = 10.36385678204158
$values[0,10] = 17.22428064879692
$values[0,11] = 18.30318065158757
$values[0,12] = 25.01590046532404
$values[1,0] = 10.34353185286267
$values[1,1] = 9.072090764205953
$values[1,2] = 13.6987164896094
$values[1,3] = 33.20697588229252
$values[1,4] = 34.94937441229736
$values[1,5] = 15.9605431290671
$values[1,6] = 25.9844134371041
$values[1,7] = 10.22362087341461
$values[1,8] = 13.39696777823846
$values[1,9] = 10.38344423138485
$values[1,10] = 17.09520849998974
$values[1,11] = 18.3325065530815
$values[1,12] = 25.09512274873964
$values[2,0] = 10.33232999700479
$values[2,1] = 9.065379425366974
$values[2,2] = 13.71298460085397
$values[2,3] = 33.26152811103319
$values[2,4] = 35.01453781551637
$values[2,5] = 15.99248696441017
$values[2,6] = 26.03158469013912
$values[2,7] = 10.23696847608208
$values[2,8] = 13.17344684886921
$values[2,9] = 10.39625756345143
$values[2,10] = 17.01724446427547
$values[2,11] = 18.35218613540148
$values[2,12] = 25.14855584794082
$values[3,0] = 10.32817375117347
$values[3,1] = 9.062935115917341
$values[3,2] = 13.71921253737116
$values[3,3] = 33.28550184703148
$values[3,4] = 35.043464317481
$values[3,5] = 16.00609222307081
$values[3,6] = 26.05204189248587
$values[3,7] = 10.2425999624393
$values[3,8] = 13.08147971793086
$values[3,9] = 10.40167733320299
$values[3,10] = 16.98582195271599
$values[3,11] = 18.36062727849082
$values[3,12] = 25.17153338972378
$values[4,0] = 10.32750838848442
$values[4,1] = 9.062546855569138
$values[4,2] = 13.7202716743163
$values[4,3] = 33.28958787802767
$values[4,4] = 35.04841057386863
$values[4,5] = 16.00838687744338
$values[4,6] = 26.05551331795153
$values[4,7] = 10.24354669005479
$values[4,8] = 13.06615923190087
$values[4,9] = 10.40258926776656
$values[4,10] = 16.98062603147338
$values[4,11] = 18.3620544073878
$values[4,12] = 25.1754214153131
$values[5,0] = 10.33227228542646
$values[5,1] = 9.065345280991275
$values[5,2] = 13.71306691780767
$values[5,3] = 33.2618443743231
$values[5,4] = 35.01491833477233
$values[5,5] = 15.99266806914997
$values[5,6] = 26.031855586346
$values[5,7] = 10.23704364533912
$values[5,8] = 13.17220994394744
$values[5,9] = 10.39632985312849
$values[5,10] = 17.01681924590861
$values[5,11] = 18.35229826772941
$values[5,12] = 25.14886086166682
$values[6,0] = 10.35687184251293
$values[6,1] = 9.080239473074981
$values[6,2] = 13.68490329377819
$values[6,3] = 33.15472947200929
$values[6,4] = 34.88796711553132
$values[6,5] = 15.92846515120526
$values[6,6] = 25.93830089882513
$values[6,7] = 10.21003353631025
$values[6,8] = 13.63212360882655
$values[6,9] = 10.37044755272602
$values[6,10] = 17.17952459257375
$values[6,11] = 18.31294538015815
$values[6,12] = 25.04222123626507
$values[7,0] = 10.41771890320078
$values[7,1] = 9.118441344128602
$values[7,2] = 13.64268467402059
$values[7,3] = 32.99986307122725
$values[7,4] = 34.71438194232579
$values[7,5] = 15.8211016021672
$values[7,6] = 25.79389042322718
$values[7,7] = 10.16309754750255
$values[7,8] = 14.4972368402414
$values[7,9] = 10.32591334675115
$values[7,10] = 17.5076843775683
$values[7,11] = 18.24901813899956
$values[7,12] = 24.87117251597352
$values[8,0] = 10.4698627536603
$values[8,1] = 9.151843406569265
$values[8,2] = 13.6195657878797
$values[8,3] = 32.91984661614007
$values[8,4] = 34.63315657362287
$values[8,5] = 15.75351300330886
$values[8,6] = 25.71166667123306
$values[8,7] = 10.13226159909003
$values[8,8] = 15.10417558988707
$values[8,9] = 10.29695877723942
$values[8,10] = 17.75281069217648
$values[8,11] = 18.21007952623941
$values[8,12] = 24.76879462820184
$values[9,0] = 10.49514151456395
$values[9,1] = 9.168162740151544
$values[9,2] = 13.61075612768409
$values[9,3] = 32.89079571713805
$values[9,4] = 34.60632248389429
$values[9,5] = 15.72521495700417
$values[9,6] = 25.6794571853684
$values[9,7] = 10.11901959855367
$values[9,8] = 15.37290606550009
$values[9,9] = 10.28459821104811
$values[9,10] = 17.86487173712717
$values[9,11] = 18.19409923110705
$values[9,12] = 24.72729414194041
$values[10,0] = 10.50493262381697
$values[10,1] = 9.174500965758501
$values[10,2] = 13.60766494210441
$values[10,3] = 32.88085259546004
$values[10,4] = 34.59761947519446
$values[10,5] = 15.71485117074796
$values[10,6] = 25.66800800296306
$values[10,7] = 10.1141176857656
$values[10,8] = 15.47352188527023
$values[10,9] = 10.28003376260295
$values[10,10] = 17.90735794938544
$values[10,11] = 18.18829635598145
$values[10,12] = 24.71230954321607
$values[11,0] = 10.50281430140066
$values[11,1] = 9.17312892192861
$values[11,2] = 13.60831980573066
$values[11,3] = 32.88294696516673
$values[11,4] = 34.59942890559144
$values[11,5] = 15.71706753916919
$values[11,6] = 25.67044051600512
$values[11,7] = 10.11516840119476
$values[11,8] = 15.45190476354531
$values[11,9] = 10.28101163511488
$values[11,10] = 17.89820595810774
$values[11,11] = 18.18953506728115
$values[11,12] = 24.7155042232946
$values[12,0] = 10.49594268644575
$values[12,1] = 9.168681033449598
$values[12,2] = 13.61049691090406
$values[12,3] = 32.88995648558032
$values[12,4] = 34.60557723309415
$values[12,5] = 15.72435526607395
$values[12,6] = 25.67850026169898
$values[12,7] = 10.11861406192456
$values[12,8] = 15.38120723876273
$values[12,9] = 10.28422036394836
$values[12,10] = 17.86836624013119
$values[12,11] = 18.19361684845148
$values[12,12] = 24.72604669712558
$values[13,0] = 10.49176192374032
$values[13,1] = 9.165977108655241
$values[13,2] = 13.61186231609007
$values[13,3] = 32.89438780341464
$values[13,4] = 34.60953329758625
$values[13,5] = 15.72886505908332
$values[13,6] = 25.68353450443988
$values[13,7] = 10.12073927361229
$values[13,8] = 15.33775119906953
$values[13,9] = 10.28620092794989
$values[13,10] = 17.85009437295884
$values[13,11] = 18.19614940064527
$values[13,12] = 24.73259947077253
$values[14,0] = 10.46824162040785
$values[14,1] = 9.150799253259409
$values[14,2] = 13.6201758186216
$values[14,3] = 32.92189313967062
$values[14,4] = 34.63511414766594
$values[14,5] = 15.75541161426176
$values[14,6] = 25.7138762432122
$values[14,7] = 10.13314276564222
$values[14,8] = 15.08645752769094
$values[14,9] = 10.29778285528589
$values[14,10] = 17.7454959295417
$values[14,11] = 18.21115869032516
$values[14,12] = 24.77160896270736
$values[15,0] = 10.45420796413464
$values[15,1] = 9.141773883582536
$values[15,2] = 13.62571269280552
$values[15,3] = 32.94064989442336
$values[15,4] = 34.65340106493903
$values[15,5] = 15.77232414721887
$values[15,6] = 25.73382103813397
$values[15,7] = 10.14095278835601
$values[15,8] = 14.93034526511799
$values[15,9] = 10.30509543163133
$values[15,10] = 17.68144941506159
$values[15,11] = 18.22080978872636
$values[15,12] = 24.79684005709277
$values[16,0] = 10.44628319222086
$values[16,1] = 9.136688783903249
$values[16,2] = 13.62905809147313
$values[16,3] = 32.95212996704461
$values[16,4] = 34.66487118098207
$values[16,5] = 15.78228223885516
$values[16,6] = 25.74578164571096
$values[16,7] = 10.1455188584837
$values[16,8] = 14.83986470296202
$values[16,9] = 10.30937779175057
$values[16,10] = 17.6446651857665
$values[16,11] = 18.2265240122354
$values[16,12] = 24.81182961055299
$values[17,0] = 10.44362541364
$values[17,1] = 9.13498536850056
$values[17,2] = 13.63021841174202
$values[17,3] = 32.95613568133409
$values[17,4] = 34.66891813473086
$values[17,5] = 15.7856934575499
$values[17,6] = 25.74991523944282
$values[17,7] = 10.14707756483782
$values[17,8] = 14.80911399486523
$values[17,9] = 10.31084085217624
$values[17,10] = 17.63222075013675
$values[17,11] = 18.22848679549409
$values[17,12] = 24.81698674681065
$values[18,0] = 10.45568669412088
$values[18,1] = 9.142723696562767
$values[18,2] = 13.62510665225489
$values[18,3] = 32.93858161038364
$values[18,4] = 34.65135583595536
$values[18,5] = 15.77049992848917
$values[18,6] = 25.73164727554352
$values[18,7] = 10.14011374728674
$values[18,8] = 14.94703563696488
$values[18,9] = 10.30430909532102
$values[18,10] = 17.68826192985545
$values[18,11] = 18.21976553071208
$values[18,12] = 24.79410475512455
$values[19,0] = 10.49795515880547
$values[19,1] = 9.169983212594763
$values[19,2] = 13.60985080284972
$values[19,3] = 32.88786890410843
$values[19,4] = 34.60373171121746
$values[19,5] = 15.72220512895987
$values[19,6] = 25.67611261630522
$values[19,7] = 10.11759893654389
$values[19,8] = 15.40200458727408
$values[19,9] = 10.28327473118154
$values[19,10] = 17.8771297256337
$values[19,11] = 18.19241119129612
$values[19,12] = 24.72293027416022
$values[20,0] = 10.52685143906266
$values[20,1] = 9.188720825628923
$values[20,2] = 13.60130701359436
$values[20,3] = 32.86089122773584
$values[20,4] = 34.5811088888284
$values[20,5] = 15.69269381610073
$values[20,6] = 25.64417684800319
$values[20,7] = 10.10354002548465
$values[20,8] = 15.69263509953569
$values[20,9] = 10.27020485423625
$values[20,10] = 18.00085073823206
$values[20,11] = 18.1759817653106
$values[20,12] = 24.68067351426538
$values[21,0] = 10.51131449494369
$values[21,1] = 9.178636965819688
$values[21,2] = 13.60573667089935
$values[21,3] = 32.87472528746953
$values[21,4] = 34.59240410875231
$values[21,5] = 15.70825678950844
$values[21,6] = 25.6608224381079
$values[21,7] = 10.1109836537604
$values[21,8] = 15.53816153010405
$values[21,9] = 10.27711864978084
$values[21,10] = 17.93480172877297
$values[21,11] = 18.18461817808389
$values[21,12] = 24.70283651566605
$values[22,0] = 10.45501771333049
$values[22,1] = 9.142293962847212
$values[22,2] = 13.62538013798441
$values[22,3] = 32.93951451194062
$values[22,4] = 34.65227750416214
$values[22,5] = 15.77132392622728
$values[22,6] = 25.73262849466066
$values[22,7] = 10.14049284104111
$values[22,8] = 14.93949218552087
$values[22,9] = 10.30466435415805
$values[22,10] = 17.68518187461548
$values[22,11] = 18.22023712362964
$values[22,12] = 24.7953398780564
$values[23,0] = 10.39993273317789
$values[23,1] = 9.107158483970927
$values[23,2] = 13.65271628516495
$values[23,3] = 33.03583665674864
$values[23,4] = 34.75323250134138
$values[23,5] = 15.84816267956793
$values[23,6] = 25.82876907294693
$values[23,7] = 10.17515231391428
$values[23,8] = 14.26781356500503
$values[23,9] = 10.32591334675115
$values[23,10] = 17.5076843775683
$values[23,11] = 18.24901813899956
$values[23,12] = 24.87117251597352

$ws.Range("C2:O25").Value = $values
